# Apply crypto price/volume updates to match target OOXML diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "29.452.81"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.851.34"
$ws.Range("E3").Value = "  +0.51%  "
Set-CellText "D4" "1.0000"
$ws.Range("E4").Value = "  +0.09%  "
Set-CellText "D5" "240.80"
$ws.Range("E5").Value = "  +0.68%  "
Set-CellText "D6" "0.6304"
$ws.Range("E6").Value = "  +0.33%  "
Set-CellText "D7" "1.000"
$ws.Range("E7").Value = "  +0.03%  "
Set-CellText "D8" "0.07699"
$ws.Range("E8").Value = "  +2.38%  "
Set-CellText "D9" "0.2931"
$ws.Range("E9").Value = "  -0.34%  "
Set-CellText "D10" "24.75"
$ws.Range("E10").Value = "  +0.83%  "
Set-CellText "D11" "0.07741"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "1.883.38"
$ws.Range("E12").Value = "  +2.52%  "
Set-CellText "D13" "5.037"
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-CellText "D14" "0.00001076"
$ws.Range("E14").Value = "  +5.39%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-CellText "D15" "0.6800"
$ws.Range("E15").Value = "  +0.45%  "
Set-CellText "D16" "83.74"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "2.125.90"
$ws.Range("E17").Value = "  +2.13%  "
Set-CellText "D18" "6.199"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "29.472.95"
$ws.Range("E19").Value = "  +0.44%  "
Set-CellText "D20" "228.73"
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("E22").Value = "  +0.01%  "
Set-CellText "D23" "7.459"
$ws.Range("E23").Value = "  +0.32%  "
Set-CellText "D24" "1.000"
$ws.Range("E24").Value = "  +0.02%  "
Set-CellText "D25" "157.44"
$ws.Range("E25").Value = "  +0.60%  "
Set-CellText "D26" "0.1384"
$ws.Range("E26").Value = "  -0.39%  "
Set-CellText "D27" "8.413"
$ws.Range("E27").Value = "  +0.87%  "
Set-CellText "D28" "17.70"
$ws.Range("E28").Value = "  +0.74%  "
Set-CellText "D29" "1.333"
$ws.Range("E29").Value = "  +5.24%  "
Set-CellText "D30" "1.468"
Set-CellText "D31" "0.05680"
$ws.Range("E31").Value = "  +0.67%  "
Set-CellText "D32" "4.129"
$ws.Range("E32").Value = "  +0.43%  "
Set-CellText "D33" "4.049"
$ws.Range("E33").Value = "  +0.62%  "
Set-CellText "D34" "1.853"
$ws.Range("E34").Value = "  +1.26%  "
Set-CellText "D35" "1.165"
$ws.Range("E35").Value = "  +1.05%  "
Set-CellText "D36" "0.7081"
$ws.Range("E36").Value = "  -0.28%  "
Set-CellText "D37" "2.586"
$ws.Range("E37").Value = "  -0.27%  "
Set-CellText "D38" "2.783"
$ws.Range("E38").Value = "  +0.45%  "
Set-CellText "D39" "0.01792"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "1.219.95"
$ws.Range("E40").Value = "  -1.62%  "
Set-CellText "D41" "6.546"
$ws.Range("E41").Value = "  +5.34%  "
Set-CellText "D42" "0.9074"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-CellText "D43" "1.001"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-CellText "D44" "101.82"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-CellText "D45" "66.52"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-CellText "D46" "0.00000000119"
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-CellText "D47" "7.134"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-CellText "D48" "0.4024"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText "D49" "9.011"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText "D50" "1.688"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-CellText "D51" "0.1144"
$ws.Range("E51").Value = "  +2.56%  "
